# Add 2022-Q1 data.
#
# Before:  ... 2021-Q4, 总计
# After:   ... 2021-Q4, 2022-Q1, 总计
#
# The existing "总计" sheet (sheetId 6) is renamed to "2022-Q1" and filled
# with the new quarter's per-fund holdings; a brand-new "总计" sheet
# (taking the next sheetId) is appended with the refreshed summary table
# (old rows shifted down by one, plus a new first row for 2022-Q1).

$wb = $excel.ActiveWorkbook

# A sheet we never touch, used purely as a formatting donor (style index 2:
# centered/bordered header + index-column look used throughout this workbook).
$fmtSrc = $wb.Worksheets.Item(5)

# ---------------------------------------------------------------------
# Step 1: turn the old "总计" sheet into the new "2022-Q1" fund table.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(6)
$q1.Name = "2022-Q1"

# Wipe whatever was there before.
$q1.Cells.Clear()

# Header row.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$fundRows = @(
    @("512980", "广发中证传媒ETF", "44.11", "99.38", "3.43", "1.5130", 6),
    @("159869", "华夏中证动漫游戏ETF", "6.20", "98.75", "7.08", "0.4390", 3),
    @("516010", "国泰中证动漫游戏ETF", "4.95", "98.91", "6.95", "0.3440", 5),
    @("160629", "鹏华中证传媒指数（LOF）", "7.63", "92.90", "3.09", "0.2358", 10),
    @("161030", "富国中证体育产业指数", "2.32", "93.75", "4.73", "0.1097", 6),
    @("516770", "华泰柏瑞中证动漫游戏ETF", "1.11", "96.56", "6.90", "0.0766", 3),
    @("164818", "工银瑞信中证传媒指数（LOF）A", "1.99", "92.70", "3.19", "0.0635", 6),
    @("159805", "鹏华中证传媒ETF", "1.73", "96.29", "3.24", "0.0561", 8),
    @("159725", "工银瑞信中证线上消费主题交易型开放式指数证券投资基金", "0.75", "98.18", "2.65", "0.0199", 10),
    @("003397", "银华体育文化灵活配置混合", "0.39", "83.61", "3.76", "0.0147", 8),
    @("010677", "工银瑞信中证传媒指数（LOF）C", "0.25", "92.70", "3.19", "0.0080", 6),
    @("005965", "安信中证500指数增强A", "0.37", "89.61", "0.85", "0.0031", 9),
    @("516190", "华夏中证文娱传媒ETF", "0.07", "96.81", "2.81", "0.0020", 10),
    @("005966", "安信中证500指数增强C", "0.17", "89.61", "0.85", "0.0014", 9)
)

$r = 2
foreach ($row in $fundRows) {
    $q1.Range("A" + $r).Value = ($r - 2)
    $q1.Range("B" + $r).Value = "'" + $row[0]
    $q1.Range("C" + $r).Value = "'" + $row[1]
    $q1.Range("D" + $r).Value = "'" + $row[2]
    $q1.Range("E" + $r).Value = "'" + $row[3]
    $q1.Range("F" + $r).Value = "'" + $row[4]
    $q1.Range("G" + $r).Value = "'" + $row[5]
    $q1.Range("H" + $r).Value = $row[6]
    $r = $r + 1
}
$lastFundRow = $r - 1

# Apply the shared "index/header" formatting (style 2) to the header row and
# the numbered index column, same as every other quarterly sheet.
$fmtSrc.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$fmtSrc.Range("A2:A" + $lastFundRow).Copy()
$q1.Range("A2:A" + $lastFundRow).PasteSpecial(-4122)

$q1.Range("A1").Select()

# ---------------------------------------------------------------------
# Step 2: append a brand-new "总计" sheet with the refreshed summary.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$summaryRows = @(
    @("2022-Q1", 14, 2.89),
    @("2021-Q4", 18, 3.4),
    @("2021-Q3", 17, 2.98),
    @("2021-Q2", 6, 0.57),
    @("2021-Q1", 10, 2.18),
    @("2020-Q4", 6, 1.29)
)

$r = 2
foreach ($row in $summaryRows) {
    $total.Range("A" + $r).Value = ($r - 2)
    $total.Range("B" + $r).Value = "'" + $row[0]
    $total.Range("C" + $r).Value = $row[1]
    $total.Range("D" + $r).Value = $row[2]
    $r = $r + 1
}
$lastTotalRow = $r - 1

$fmtSrc.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$fmtSrc.Range("A2:A" + $lastTotalRow).Copy()
$total.Range("A2:A" + $lastTotalRow).PasteSpecial(-4122)

$total.Range("A1").Select()
